$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 417.8
$ws.Range("I32").Value = 250.5
$ws.Range("J32").Value = 529.3333
$ws.Range("K32").Value = 250.5
$ws.Range("L32").Value = 529.3333
$ws.Range("M32").Value = 75.5
$ws.Range("N32").Value = -1181.3333
$ws.Range("H33").Value = 370
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H70").Value = 1366.6666
$ws.Range("I70").Value = 1466.6666
$ws.Range("J70").Value = 1166.6666
$ws.Range("K70").Value = 4399.9998
$ws.Range("L70").Value = 3499.9998
$ws.Range("M70").Value = -4129.9998
$ws.Range("N70").Value = -4039.9998
$ws.Range("H73").Value = 1366.6666
$ws.Range("I73").Value = 1466.6666
$ws.Range("J73").Value = 1166.6666
$ws.Range("K73").Value = 4399.9998
$ws.Range("L73").Value = 3499.9998
$ws.Range("M73").Value = -3463.9998
$ws.Range("N73").Value = -5371.9998
$ws.Range("H98").Value = 392.3846
$ws.Range("I98").Value = 392.3846
$ws.Range("K98").Value = 392.3846
$ws.Range("M98").Value = 1105.6154
$ws.Range("H100").Value = 3534.5386
$ws.Range("I100").Value = 3333.2222
$ws.Range("J100").Value = 3987.5
$ws.Range("K100").Value = 3333.2222
$ws.Range("L100").Value = 3987.5
$ws.Range("M100").Value = -2792.2222
$ws.Range("N100").Value = -5069.5
$ws.Range("H122").Value = 392.3846
$ws.Range("I122").Value = 392.3846
$ws.Range("K122").Value = 1177.1538
$ws.Range("M122").Value = 1272.8462
$ws.Range("H129").Value = 1310.7
$ws.Range("J129").Value = 1324.4407
$ws.Range("L129").Value = 3973.3221
$ws.Range("N129").Value = -13973.3221

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1715.5238
$ws.Range("I2").Value = 1286.7142
$ws.Range("K2").Value = 1286.7142
$ws.Range("M2").Value = -1173.7142
$ws.Range("H102").Value = 1609.909
$ws.Range("I102").Value = 1387.1428
$ws.Range("K102").Value = 1387.1428
$ws.Range("M102").Value = 234.8571999999999
$ws.Range("H116").Value = 1715.5238
$ws.Range("I116").Value = 1286.7142
$ws.Range("K116").Value = 1286.7142
$ws.Range("M116").Value = 1007.2858
$ws.Range("H122").Value = 1944
$ws.Range("I122").Value = 1920.2667
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 5760.800099999999
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -3310.800099999999
$ws.Range("N122").Value = -11800

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1715.5238
$ws.Range("I3").Value = 1286.7142
$ws.Range("K3").Value = 1286.7142
$ws.Range("M3").Value = -1172.7142
$ws.Range("H27").Value = 42500
$ws.Range("J27").Value = 42500
$ws.Range("L27").Value = 42500
$ws.Range("N27").Value = -42884
$ws.Range("H86").Value = 1412.7819
$ws.Range("I86").Value = 1341.4318
$ws.Range("J86").Value = 1698.1818
$ws.Range("K86").Value = 1341.4318
$ws.Range("L86").Value = 1698.1818
$ws.Range("M86").Value = -218.4318000000001
$ws.Range("N86").Value = -3944.1818
$ws.Range("H89").Value = 1412.7819
$ws.Range("I89").Value = 1341.4318
$ws.Range("J89").Value = 1698.1818
$ws.Range("K89").Value = 6707.159000000001
$ws.Range("L89").Value = 8490.909
$ws.Range("M89").Value = -1091.159000000001
$ws.Range("N89").Value = -19722.909

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 10370.333
$ws.Range("I50").Value = 4000
$ws.Range("J50").Value = 13555.5
$ws.Range("K50").Value = 4000
$ws.Range("L50").Value = 13555.5
$ws.Range("M50").Value = -3375
$ws.Range("N50").Value = -14805.5
$ws.Range("H51").Value = 10845
$ws.Range("I51").Value = 6690
$ws.Range("K51").Value = 6690
$ws.Range("M51").Value = -5954
$ws.Range("H60").Value = 11638.286
$ws.Range("J60").Value = 11970.2
$ws.Range("L60").Value = 11970.2
$ws.Range("N60").Value = -12992.2
$ws.Range("H61").Value = 10845
$ws.Range("I61").Value = 6690
$ws.Range("K61").Value = 6690
$ws.Range("M61").Value = -6342
$ws.Range("H99").Value = 20004498
$ws.Range("I99").Value = 3479.4666
$ws.Range("J99").Value = 50006028
$ws.Range("K99").Value = 3479.4666
$ws.Range("L99").Value = 50006028
$ws.Range("M99").Value = -1981.4666
$ws.Range("N99").Value = -50009024
$ws.Range("H122").Value = 1453.762
$ws.Range("I122").Value = 1367.5555
$ws.Range("K122").Value = 4102.666499999999
$ws.Range("M122").Value = -1652.666499999999
$ws.Range("H126").Value = 20004498
$ws.Range("I126").Value = 3479.4666
$ws.Range("J126").Value = 50006028
$ws.Range("K126").Value = 10438.3998
$ws.Range("L126").Value = 150018084
$ws.Range("M126").Value = -7968.399800000001
$ws.Range("N126").Value = -150023024
$ws.Range("H132").Value = 56312.5
$ws.Range("I132").Value = 88018.836
$ws.Range("J132").Value = 8753
$ws.Range("K132").Value = 264056.508
$ws.Range("L132").Value = 26259
$ws.Range("M132").Value = -261526.508
$ws.Range("N132").Value = -31319
$ws.Range("H135").Value = 35372.668
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 35372.668
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 35372.668
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -45512.668

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1911.4529
$ws.Range("I2").Value = 2295.932
$ws.Range("J2").Value = 31.777779
$ws.Range("K2").Value = 13775.592
$ws.Range("L2").Value = 190.666674
$ws.Range("M2").Value = -13662.592
$ws.Range("N2").Value = -416.666674
$ws.Range("H68").Value = 50168
$ws.Range("I68").Value = 333
$ws.Range("J68").Value = 100003
$ws.Range("K68").Value = 999
$ws.Range("L68").Value = 300009
$ws.Range("M68").Value = -188
$ws.Range("N68").Value = -301631
$ws.Range("H71").Value = 50168
$ws.Range("I71").Value = 333
$ws.Range("J71").Value = 100003
$ws.Range("K71").Value = 2997
$ws.Range("L71").Value = 900027
$ws.Range("M71").Value = 1059
$ws.Range("N71").Value = -908139
$ws.Range("H92").Value = 599.8889
$ws.Range("I92").Value = 666.6667
$ws.Range("J92").Value = 466.33334
$ws.Range("K92").Value = 2000.0001
$ws.Range("L92").Value = 1399.00002
$ws.Range("M92").Value = -752.0001
$ws.Range("N92").Value = -3895.00002
$ws.Range("H107").Value = 12677.5
$ws.Range("J107").Value = 223.4
$ws.Range("L107").Value = 670.2
$ws.Range("N107").Value = -4510.2
$ws.Range("H122").Value = 372.53845
$ws.Range("I122").Value = 314
$ws.Range("J122").Value = 409.125
$ws.Range("K122").Value = 2826
$ws.Range("L122").Value = 3682.125
$ws.Range("M122").Value = -376
$ws.Range("N122").Value = -8582.125
$ws.Range("H131").Value = 170322.55
$ws.Range("I131").Value = 952
$ws.Range("J131").Value = 186005
$ws.Range("K131").Value = 2856
$ws.Range("L131").Value = 558015
$ws.Range("M131").Value = 2184
$ws.Range("N131").Value = -568095

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3556.1875
$ws.Range("I80").Value = 3044.4443
$ws.Range("K80").Value = 3044.4443
$ws.Range("M80").Value = -2046.4443
$ws.Range("H83").Value = 3556.1875
$ws.Range("I83").Value = 3044.4443
$ws.Range("K83").Value = 15222.2215
$ws.Range("M83").Value = -10230.2215
$ws.Range("H122").Value = 2871.2
$ws.Range("I122").Value = 2180
$ws.Range("K122").Value = 6540
$ws.Range("M122").Value = -4090

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1592.4138
$ws.Range("I82").Value = 1851.2858
$ws.Range("J82").Value = 1350.8
$ws.Range("K82").Value = 1851.2858
$ws.Range("L82").Value = 1350.8
$ws.Range("M82").Value = -1490.2858
$ws.Range("N82").Value = -2072.8
$ws.Range("H85").Value = 1592.4138
$ws.Range("I85").Value = 1851.2858
$ws.Range("J85").Value = 1350.8
$ws.Range("K85").Value = 1851.2858
$ws.Range("L85").Value = 1350.8
$ws.Range("M85").Value = -603.2858000000001
$ws.Range("N85").Value = -3846.8
$ws.Range("H132").Value = 2540.3125
$ws.Range("I132").Value = 1863.8572
$ws.Range("J132").Value = 3066.4443
$ws.Range("K132").Value = 5591.571599999999
$ws.Range("L132").Value = 9199.332900000001
$ws.Range("M132").Value = -3061.571599999999
$ws.Range("N132").Value = -14259.3329
$ws.Range("H136").Value = 56880.555
$ws.Range("I136").Value = 63740.625
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 191221.875
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -188671.875
$ws.Range("N136").Value = -11100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 111112230
$ws.Range("I81").Value = 1348.3334
$ws.Range("J81").Value = 333334000
$ws.Range("K81").Value = 2696.6668
$ws.Range("L81").Value = 666668000
$ws.Range("M81").Value = -1635.6668
$ws.Range("N81").Value = -666670122
$ws.Range("H84").Value = 111112230
$ws.Range("I84").Value = 1348.3334
$ws.Range("J84").Value = 333334000
$ws.Range("K84").Value = 13483.334
$ws.Range("L84").Value = 3333340000
$ws.Range("M84").Value = -8179.333999999999
$ws.Range("N84").Value = -3333350608
$ws.Range("H118").Value = 43346
$ws.Range("J118").Value = 43346
$ws.Range("L118").Value = 43346
$ws.Range("N118").Value = -46660
$ws.Range("H132").Value = 1324.1333
$ws.Range("I132").Value = 758.7692
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 2276.3076
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = 253.6923999999999
$ws.Range("N132").Value = -20057
